$wb = $excel.ActiveWorkbook

$wsDone = $wb.Worksheets.Item("done")
$wsShort = $wb.Worksheets.Item("short term")

# The task that used to be row 21 on "short term" is now finished, so it
# moves to the bottom of the "done" list (new row 67), keeping its original
# text and assignee.
$doneTaskText = $wsShort.Range("A21").Value2
$doneTaskPerson = $wsShort.Range("B21").Value2

$newRow = 67

# Match the formatting used by the rest of the "done" list: column A wraps
# text and is bordered (same look as the row above it), column B uses the
# plain highlighted "person" look used throughout column B.
$wsDone.Range("A66").Copy()
$wsDone.Range("A" + $newRow).PasteSpecial(-4122) | Out-Null
$wsDone.Range("B44").Copy()
$wsDone.Range("B" + $newRow).PasteSpecial(-4122) | Out-Null
$wsDone.Application.CutCopyMode = $false

$wsDone.Range("A" + $newRow).Value2 = $doneTaskText
$wsDone.Range("B" + $newRow).Value2 = $doneTaskPerson

# This task's text wraps onto two lines at column A's width, just like it
# did on its previous sheet - match that row height.
$wsDone.Rows.Item($newRow).RowHeight = 29

# Update the sheet selection to mirror the new last row.
$wsDone.Range("A" + $newRow).Select()

# Replace the now-completed row 21 on "short term" with the next new task,
# dropping the old person/highlight formatting (delete + retype clears the
# stale explicit row height / fill that belonged to the finished task).
$wsShort.Rows.Item(21).Delete()
$wsShort.Range("A21").Value2 = "88. change the examples of DRomics data results (triclosan) to ass the column yatdosemax"
